$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format temporarily so numeric-looking strings
# (e.g. "1.001") are stored as text, matching the inlineStr cells in the source.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.997.10"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.640.55"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "209.15"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "0.5158"
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "0.2562"
$ws.Range("E8").Value = "  -2.92%  "
$ws.Range("D9").Value = "0.06210"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").Value = "20.39"
$ws.Range("E10").Value = "  -3.62%  "
$ws.Range("D11").Value = "0.07523"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "1.635.89"
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").Value = "4.354"
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("D14").Value = "1.855.43"
$ws.Range("D15").Value = "0.5382"
$ws.Range("E15").Value = "  -3.74%  "
$ws.Range("D16").Value = "0.0₅7950"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "64.99"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("D18").Value = "25.991.72"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").Value = "4.645"
$ws.Range("E20").Value = "  -2.96%  "
$ws.Range("D21").Value = "186.04"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "10.01"
$ws.Range("E22").Value = "  -3.27%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "6.080"
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").Value = "145.69"
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "7.333"
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "0.1191"
$ws.Range("E27").Value = "  -4.28%  "
$ws.Range("D28").Value = "15.47"
$ws.Range("E28").Value = "  -2.84%  "
$ws.Range("D29").Value = "1.373"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("D30").Value = "0.05951"
$ws.Range("E30").Value = "  -3.87%  "
$ws.Range("D31").Value = "1.241"
$ws.Range("E31").Value = "  -3.01%  "
$ws.Range("D32").Value = "3.355"
$ws.Range("E32").Value = "  -3.39%  "
$ws.Range("D33").Value = "3.345"
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("D34").Value = "1.604"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").Value = "0.9685"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("D36").Value = "2.372"
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("D37").Value = "2.724"
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("D38").Value = "0.5817"
$ws.Range("E38").Value = "  -3.74%  "
$ws.Range("D39").Value = "0.01583"
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("D40").Value = "1.051.69"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("D41").Value = "5.770"
$ws.Range("E41").Value = "  -5.83%  "
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").Value = "0.8382"
$ws.Range("E43").Value = "  -3.05%  "
$ws.Range("D44").Value = "99.68"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "1.783.10"
$ws.Range("E45").Value = "  -1.98%  "
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("D47").Value = "1.007"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").Value = "54.17"
$ws.Range("E48").Value = "  -3.17%  "
$ws.Range("D49").Value = "0.05198"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").Value = "7.871"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").Value = "0.4223"
$ws.Range("E51").Value = "  -0.69%  "

# Restore default (Normal) style on column D so no stray style index is left
# on cells that originally had none.
$ws.Range("D2:D51").Style = "Normal"

